# Regenerate the "K" (strikeout) column values for the save_data sheet.
# The workbook's G column (header "K") holds per-game values that were
# recomputed upstream (regen std/mean, calc and write s_vals) and now
# need to be written back into rows 2-32.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 2
    3  = 2
    4  = 1
    5  = 2
    6  = 1
    7  = 1
    8  = 2
    9  = 1
    10 = 1
    11 = 2
    12 = 2
    13 = 2
    14 = 1
    15 = 0
    16 = 0
    17 = 2
    18 = 2
    19 = 3
    20 = 1
    21 = 1
    22 = 0
    23 = 0
    24 = 1
    25 = 2
    26 = 1
    27 = 1
    28 = 0
    29 = 0
    30 = 1
    31 = 2
    32 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
